$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "2022-Q3" worksheet before the existing "2022-Q2" sheet ---
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

# Header row (row 1) - mirrors the other quarter sheets
$q3Sheet.Range("B1").Value = "基金代码"
$q3Sheet.Range("C1").Value = "基金名称"
$q3Sheet.Range("D1").Value = "基金规模"
$q3Sheet.Range("E1").Value = "股票总仓位"
$q3Sheet.Range("F1").Value = "仓位占比"
$q3Sheet.Range("G1").Value = "持有市值(亿元)"
$q3Sheet.Range("H1").Value = "仓位排名"

# Data row (row 2)
$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").Value = "166024"
$q3Sheet.Range("C2").Value = "中欧恒利三年定期开放混合"
$q3Sheet.Range("D2").Value = "3.99"
$q3Sheet.Range("E2").Value = "98.45"
$q3Sheet.Range("F2").Value = "3.69"
$q3Sheet.Range("G2").Value = "0.1472"
$q3Sheet.Range("H2").Value = 9

# Match the header/data cell styling used on the other quarter sheets
$q3Sheet.Range("B1:H1").Style = $q2Sheet.Range("B1:H1").Style
$q3Sheet.Range("A2").Style = $q2Sheet.Range("A2").Style

# --- 2. Insert a new row into "总计" summarising the new quarter ---
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.15

# Keep the "编号" column (A) sequential (0,1,2) after the insert
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
